$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2
$ws.Range("C2").Value = 3231
$ws.Range("E2").Value = "23 -> 67 -> 70 -> 74 -> 79 -> 2 -> 26 -> 8 -> 20 -> 46 -> 49 -> 23"
$ws.Range("F2").Value = "23 -> 22 -> 21 -> 67 -> 68 -> 70 -> 72 -> 74 -> 76 -> 77 -> 78 -> 79 -> 80 -> 2 -> 3 -> 25 -> 26 -> 27 -> 6 -> 5 -> 8 -> 12 -> 16 -> 17 -> 20 -> 24 -> 32 -> 44 -> 46 -> 45 -> 49 -> 45 -> 43 -> 32 -> 24 -> 23"
$ws.Range("G2").Value = 0.07203245162963867
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.07203245162963867

# Row 3
$ws.Range("C3").Value = 4082
$ws.Range("F3").Value = "48 -> 47 -> 46 -> 44 -> 41 -> 32 -> 24 -> 23 -> 19 -> 16 -> 15 -> 14 -> 10 -> 7 -> 4 -> 1 -> 77 -> 78 -> 79 -> 81 -> 3 -> 25 -> 26 -> 27 -> 28 -> 29 -> 30 -> 37 -> 38 -> 41 -> 32 -> 24 -> 23 -> 22 -> 66 -> 65 -> 64 -> 61 -> 62 -> 58 -> 57 -> 55 -> 53 -> 49 -> 45 -> 46 -> 47 -> 48"
$ws.Range("G3").Value = 0.06498193740844727
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.06498193740844727

# Row 4
$ws.Range("C4").Value = 3796
$ws.Range("E4").Value = "11 -> 63 -> 64 -> 21 -> 76 -> 1 -> 34 -> 29 -> 31 -> 40 -> 17 -> 11"
$ws.Range("F4").Value = "11 -> 10 -> 14 -> 18 -> 21 -> 22 -> 66 -> 65 -> 64 -> 61 -> 62 -> 63 -> 62 -> 61 -> 64 -> 65 -> 66 -> 22 -> 21 -> 67 -> 68 -> 70 -> 72 -> 74 -> 76 -> 77 -> 1 -> 2 -> 3 -> 25 -> 26 -> 27 -> 28 -> 34 -> 28 -> 29 -> 30 -> 31 -> 38 -> 39 -> 40 -> 39 -> 38 -> 31 -> 20 -> 17 -> 16 -> 12 -> 11"
$ws.Range("G4").Value = 0.08513379096984863
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.08513379096984863
